# Add a new "Descriptors" sheet, positioned right after "QuantityValue"
# and right before "Acquisition", with header row: descriptor_name, descriptor_thing.

$wb = $excel.ActiveWorkbook

# Create the new worksheet and name it.
$ws = $wb.Worksheets.Add()
$ws.Name = "Descriptors"

# Populate the header row before moving the sheet (writes target the
# still-active freshly-added sheet reliably at this point).
$ws.Range("A1").Value = "descriptor_name"
$ws.Range("B1").Value = "descriptor_thing"

# Reposition it immediately before the "Acquisition" sheet, i.e. right
# after "QuantityValue".
$target = $wb.Worksheets.Item("Acquisition")
$ws.Move($target)
